# Applies the change described by the diff: the single "step" row that used
# to belong to test case TC4 ("Clica para realizar o empenho de uma
# diária." / "SYSTEM Apresenta a tela de Registrar Empenho") now appears
# under TC3, and the step row that used to belong to TC3 ("Clica para
# atribuir/desatribuir o registro a si mesmo." / "SYSTEM Atualiza a lista
# de registros de solicitações...") now appears under TC4. The "Test Case
# ID" labels themselves (TC3, TC4) stay where they are; only the Steps
# (column B) and Expected Results (column D) content of their single step
# row is swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Test Case ID:" header rows for TC3 and TC4 by searching for
# their labels, so the edit does not depend on assuming fixed row numbers.
$tc3Header = $ws.Cells.Find("TC3")
$tc4Header = $ws.Cells.Find("TC4")

# Each test-case block follows the fixed layout:
#   row+0 = Test Case ID / Priority / Executed by
#   row+1 = Description
#   row+2 = Precondition
#   row+3 = # / Steps / Test Data / Expected Results / ... (column headers)
#   row+4 = first (and here, only) step row, with Steps in column B and
#           Expected Results in column D
$tc3StepRow = $tc3Header.Row + 4
$tc4StepRow = $tc4Header.Row + 4

$tc3StepCell     = $ws.Cells.Item($tc3StepRow, 2)   # Steps column (B) for TC3
$tc3ExpectedCell = $ws.Cells.Item($tc3StepRow, 4)   # Expected Results column (D) for TC3
$tc4StepCell     = $ws.Cells.Item($tc4StepRow, 2)   # Steps column (B) for TC4
$tc4ExpectedCell = $ws.Cells.Item($tc4StepRow, 4)   # Expected Results column (D) for TC4

# NOTE: use Value2 (not Value) to read/write cell contents.
$tc3StepValue     = $tc3StepCell.Value2
$tc3ExpectedValue = $tc3ExpectedCell.Value2
$tc4StepValue     = $tc4StepCell.Value2
$tc4ExpectedValue = $tc4ExpectedCell.Value2

$tc3StepCell.Value2     = $tc4StepValue
$tc3ExpectedCell.Value2 = $tc4ExpectedValue
$tc4StepCell.Value2     = $tc3StepValue
$tc4ExpectedCell.Value2 = $tc3ExpectedValue
